$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.781.75"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "2.220.66"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.08"
$ws.Range("E5").Value = "  -5.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.02"
$ws.Range("E6").Value = "  -8.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.512"
$ws.Range("E9").Value = "  -8.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.90"
$ws.Range("E10").Value = "  -9.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0794"
$ws.Range("E11").Value = "  -4.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.03"
$ws.Range("E12").Value = "  -8.42%  "
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("D14").Value = "2.557.65"
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("D15").Value = "2.252.32"
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.803"
$ws.Range("E16").Value = "  -7.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.29"
$ws.Range("E17").Value = "  -7.92%  "
$ws.Range("D18").Value = "43.549.44"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "0.0₃0943"
$ws.Range("E19").Value = "  -4.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.82"
$ws.Range("E20").Value = "  -11.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.07"
$ws.Range("E21").Value = "  -7.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "63.78"
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.59"
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.86"
$ws.Range("E24").Value = "  -8.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.90"
$ws.Range("E26").Value = "  -10.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.65"
$ws.Range("E27").Value = "  -5.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.54"
$ws.Range("E29").Value = "  -5.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.82"
$ws.Range("E30").Value = "  -6.87%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.61"
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.94"
$ws.Range("E32").Value = "  -4.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0792"
$ws.Range("E33").Value = "  -7.44%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.21"
$ws.Range("E34").Value = "  +3.13%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("E36").Value = "  -2.35%  "
$ws.Range("E37").Value = "  -9.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.73"
$ws.Range("E38").Value = "  -11.87%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.42"
$ws.Range("E39").Value = "  -10.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.71"
$ws.Range("E40").Value = "  -11.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0292"
$ws.Range("E41").Value = "  -7.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.21"
$ws.Range("E42").Value = "  -14.02%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.716.86"
$ws.Range("E44").Value = "  -4.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "83.06"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -6.78%  "
$ws.Range("E47").Value = "  -8.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.14"
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.98"
$ws.Range("E49").Value = "  -5.02%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "14.26"
$ws.Range("E50").Value = "  +3.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "66.97"
$ws.Range("E51").Value = "  -11.83%  "
